$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "74.254.81"
$ws.Range("E2").Value = "  +7.81%  "
$ws.Range("D3").Value = "2.633.64"
$ws.Range("E3").Value = "  +7.92%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'185.87"
$ws.Range("E5").Value = "  +14.37%  "
$ws.Range("D6").Value = "'582.58"
$ws.Range("E6").Value = "  +4.17%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +4.64%  "
$ws.Range("D9").Value = "'0.202"
$ws.Range("E9").Value = "  +17.41%  "
$ws.Range("D10").Value = "2.632.89"
$ws.Range("E10").Value = "  +8.00%  "
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("E12").Value = "  +8.02%  "
$ws.Range("E13").Value = "  +3.38%  "
$ws.Range("E14").Value = "  +6.16%  "
$ws.Range("D15").Value = "73.947.16"
$ws.Range("E15").Value = "  +7.57%  "
$ws.Range("D16").Value = "3.114.46"
$ws.Range("E16").Value = "  +7.78%  "
$ws.Range("D17").Value = "'26.30"
$ws.Range("E17").Value = "  +12.82%  "
$ws.Range("D18").Value = "2.635.08"
$ws.Range("E18").Value = "  +7.93%  "
$ws.Range("D19").Value = "'9.09"
$ws.Range("E19").Value = "  +30.44%  "
$ws.Range("D20").Value = "'11.82"
$ws.Range("E20").Value = "  +11.50%  "
$ws.Range("D21").Value = "'372.48"
$ws.Range("E21").Value = "  +9.70%  "
$ws.Range("E22").Value = "  +18.23%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  +4.80%  "
$ws.Range("E26").Value = "  +11.67%  "
$ws.Range("E27").Value = "  +13.76%  "
$ws.Range("D28").Value = "2.769.89"
$ws.Range("E28").Value = "  +7.90%  "
$ws.Range("E29").Value = "  +3.75%  "
$ws.Range("D30").Value = "0.0₃0951"
$ws.Range("E30").Value = "  +15.56%  "
$ws.Range("D31").Value = "'526.86"
$ws.Range("E31").Value = "  +22.41%  "
$ws.Range("E32").Value = "  +18.98%  "
$ws.Range("D33").Value = "'7.73"
$ws.Range("E33").Value = "  +8.05%  "
$ws.Range("E34").Value = "  +8.68%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'163.27"
$ws.Range("E36").Value = "  +2.64%  "
$ws.Range("E37").Value = "  +12.65%  "
$ws.Range("E38").Value = "  +6.29%  "
$ws.Range("E39").Value = "  +1.51%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  +12.19%  "
$ws.Range("E42").Value = "  +9.82%  "
$ws.Range("E43").Value = "  +10.43%  "
$ws.Range("D44").Value = "'160.60"
$ws.Range("E44").Value = "  +23.32%  "
$ws.Range("E45").Value = "  +11.75%  "
$ws.Range("E46").Value = "  +14.82%  "
$ws.Range("D47").Value = "'38.98"
$ws.Range("E47").Value = "  +3.89%  "
$ws.Range("D48").Value = "'0.0856"
$ws.Range("E48").Value = "  +18.68%  "
$ws.Range("E49").Value = "  +8.83%  "
$ws.Range("D50").Value = "'0.530"
$ws.Range("E50").Value = "  +9.97%  "
$ws.Range("D51").Value = "'21.00"
$ws.Range("E51").Value = "  +24.03%  "
